$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Alone.global.menage"
$ws.Range("C1").Value = "With the family of origin (parents, etc.).global.menage"
$ws.Range("D1").Value = "With partner / children.global.menage"
$ws.Range("E1").Value = "With friends or other people (with no family relation).global.menage"
$ws.Range("F1").Value = "In detention.global.menage"
$ws.Range("G1").Value = "In institutions /shelters (not detention).global.menage"
$ws.Range("H1").Value = "Other.global.menage"
$ws.Range("I1").Value = "Not known / missing.global.menage"
$ws.Range("J1").Value = "Total.global.menage"
